$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.785.41'
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").Value = '1.967.95'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.90'
$ws.Range("E5").Value = '  +0.35%  '

$ws.Range("E6").Value = '  +0.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.35'
$ws.Range("E7").Value = '  +1.84%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  +1.46%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0815'
$ws.Range("E10").Value = '  -3.08%  '

$ws.Range("E11").Value = '  -0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.49'
$ws.Range("E12").Value = '  +3.35%  '

$ws.Range("D13").Value = '2.259.06'
$ws.Range("E13").Value = '  +1.24%  '

$ws.Range("E14").Value = '  +0.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.79'
$ws.Range("E15").Value = '  +1.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.29'
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").Value = '1.978.30'
$ws.Range("E17").Value = '  +2.32%  '

$ws.Range("D18").Value = '36.694.22'
$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.92'
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = '0.0₃0864'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.50'
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.45'
$ws.Range("E24").Value = '  -0.09%  '

$ws.Range("E25").Value = '  +2.80%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.36'
$ws.Range("E26").Value = '  +0.69%  '

$ws.Range("E27").Value = '  +14.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.69'
$ws.Range("E28").Value = '  -1.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.44'
$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  +1.26%  '

$ws.Range("E31").Value = '  -1.53%  '

$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("E33").Value = '  -1.68%  '

$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.27'
$ws.Range("E35").Value = '  +5.69%  '

$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.10'
$ws.Range("E37").Value = '  -2.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.40'
$ws.Range("E38").Value = '  +11.17%  '

$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("E40").Value = '  +3.37%  '

$ws.Range("E41").Value = '  -2.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0213'
$ws.Range("E42").Value = '  +1.82%  '

$ws.Range("E43").Value = '  -0.74%  '

$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").Value = '1.365.76'
$ws.Range("E45").Value = '  +0.66%  '

$ws.Range("E46").Value = '  +0.84%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.08'
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.18'
$ws.Range("E48").Value = '  +0.03%  '

$ws.Range("E49").Value = '  +0.75%  '

$ws.Range("D50").Value = '2.149.47'
$ws.Range("E50").Value = '  +1.28%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.89'
$ws.Range("E51").Value = '  -3.06%  '
